$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1.045389290513603
$ws.Range("G2").Value = 0.6629871557899286
$ws.Range("I2").Value = 0.9275548951266984
$ws.Range("J2").Value = 3.323323020455365
$ws.Range("K2").Value = 0.6758592141011408
